# Edit: Rachel Lowe CE301 Poster HECC-IT
#
# 1. In "Rectangle 7" (the "WHAT IS HECC-IT?" poster panel), change the
#    wording of the final paragraph from
#       "And then, once that's done, just open HECC-UP and you have your
#        playable hypertext game!"
#    to
#       "And then, once that's done, just put it through HECC-UP and you
#        have your playable hypertext game!"
#    and append two new paragraphs (separated by blank lines) describing
#    HECC-IT as a standalone desktop app and the HTML/JS output format.
#
# 2. In "Rectangle 10" (the poster panel that held the sample .hecc story
#    script), delete all of the sample-story text, leaving an empty
#    paragraph - mirroring the now-removed "Rectangle 9" example panel.

$rsquo = [char]0x2019   # U+2019 RIGHT SINGLE QUOTATION MARK, used by the
                         # deck's body copy in place of a plain apostrophe.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        if ($slide.Shapes.Item($i).Name -eq $name) { return $slide.Shapes.Item($i) }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) "Rectangle 7" - update closing sentence + add two new paragraphs
# ---------------------------------------------------------------------
$rect7 = Get-ShapeByName $s "Rectangle 7"
$tr7 = $rect7.TextFrame.TextRange

# --- 1a. "...just open HECC-UP..." -> "...just put it through HECC-UP..."
$fullText = $tr7.Text
$oldSentence = "And then, once that" + $rsquo + "s done, just open "
$newSentence = "And then, once that" + $rsquo + "s done, just put it through "

$idx = $fullText.IndexOf("And then, once that")
if ($idx -ge 0) {
    $target = $tr7.Characters($idx + 1, $oldSentence.Length)
    $target.Text = $newSentence
}

# --- 1b. Append two new paragraphs (with a blank line before each).
#     Insert all of the plain text first (so nothing downstream inherits
#     bold formatting), then bold the two "HECC-IT" mentions afterwards.
$tr7 = $rect7.TextFrame.TextRange
$lenBeforeAppend = $tr7.Length

$para1 = "HECC-IT is a standalone desktop application, so users only need to download it and run it."
$para2 = "Games made with HECC-IT are in a client-side HTML/JavaScript format, so anyone with a web browser can download and play them."
$appendText = "`r`r" + $para1 + "`r`r" + $para2

$tr7.InsertAfter($appendText) | Out-Null
$tr7 = $rect7.TextFrame.TextRange

# Bold the "HECC-IT" at the start of para1.
$off1 = $lenBeforeAppend + 2 + 1   # skip the two blank-paragraph marks
$tr7.Characters($off1, 7).Font.Bold = $true

# Bold the "HECC-IT" that follows "Games made with " in para2.
$gamesPrefix = "Games made with "
$idxPara2InAppend = $appendText.IndexOf($para2)
$off2 = $lenBeforeAppend + $idxPara2InAppend + $gamesPrefix.Length + 1
$tr7.Characters($off2, 7).Font.Bold = $true

# ---------------------------------------------------------------------
# 2) "Rectangle 10" - clear the sample story text
# ---------------------------------------------------------------------
$rect10 = Get-ShapeByName $s "Rectangle 10"
$tr10 = $rect10.TextFrame.TextRange
$tr10.Text = ""
